# "Added to majors folder" - adds notes + status colors for several
# majors (Studio Art, Asian Studies, Chinese Studies, Japanese Studies,
# Biology, Business Administration, Chemistry) on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Status (column B) color updates -------------------------------
# Reuse existing formatting already present on the sheet so no new
# style entries are minted: B4/B5 already carry the "yellow" status
# style, B3 already carries the "green" status style.
$ws.Range("B5").Copy()
$ws.Range("B6:B7").PasteSpecial(-4122)

$ws.Range("B3").Copy()
$ws.Range("B8:B12").PasteSpecial(-4122)

# --- New notes (columns C-F) ----------------------------------------
$ws.Range("C7").Value = "`"at least three courses in a single Asian language at Furman or three non-language Asian Studies courses`" I just put all these things as requirements. We can go back later and or the requirements."
$ws.Range("C6").Value = "12 Additional credits of studio arts (where some THA classes count too)"
$ws.Range("C8").Value = "May have a requirement satisfied by an FYW"

$ws.Range("C10").Value = "MTH - 120 and 145 or MTH 151"
$ws.Range("D10").Value = "Has substitutions, including 2 mayX courses subbing for 1 elective"
$ws.Range("E10").Value = "Includes a second major for education oriented students"
$ws.Range("F10").Value = "The third elective may be satisfied by 1 couse or by 341 and 342 together"

$ws.Range("C11").Value = "Block - I just put the 4 classes as required."
$ws.Range("C12").Value = "Includes education major"

# --- Row heights ------------------------------------------------------
$ws.Range("1:66").RowHeight = 16
$ws.Range("48:48").RowHeight = 15.75
$ws.Range("62:62").RowHeight = 15.75

# --- Selection / view --------------------------------------------------
$ws.Range("A6").Select()
